# ADD results from server
# Update the row-2 numeric values ("inv cost" results) on sheets 2025, 2030,
# 2035, 2040, and 2045 to reflect freshly computed server results. Sheet 2050
# is left untouched (no changes in the source diff).

$wb = $excel.ActiveWorkbook

# Map of sheet name -> column letter -> new value (only cells that change)
$updates = @{
    "2025" = @{
        "B2" = 10372.65132737054
        "E2" = 289260.5393052954
        "G2" = 80959.25712661834
        "I2" = 161710.6685703679
        "L2" = 484922.2142001599
        "M2" = 105953.7713982
        "N2" = 70003.73489578845
        "O2" = 69744.89343456978
    }
    "2030" = @{
        "A2" = 0
        "B2" = 31203.23858116339
        "E2" = 170658.5511254234
        "I2" = 209080.6134235085
        "L2" = 63518.11613148725
        "M2" = 68536.72857011756
        "N2" = 19285.19160463996
        "O2" = 27033.1386905727
    }
    "2035" = @{
        "A2" = 27543.1755456332
        "B2" = 22113.21643273498
        "E2" = 114655.4402706629
        "I2" = 153866.0861464091
        "M2" = 44638.22942194272
        "N2" = 39676.88529639924
        "O2" = 31311.04369977792
    }
    "2040" = @{
        "N2" = 1142.580190039942
        "O2" = 0
    }
    "2045" = @{
        "A2" = 29588.33508286276
        "N2" = 4347.543515635315
        "O2" = 20429.76977394434
    }
}

foreach ($sheetName in $updates.Keys) {
    # Cast to [string] explicitly: hashtable keys that look numeric (e.g.
    # "2025") get auto-converted to Int64 by the PowerShell enumerator,
    # which would make Worksheets.Item() treat them as a 1-based index
    # instead of a sheet name.
    $sheetNameStr = [string]$sheetName
    $ws = $wb.Worksheets.Item($sheetNameStr)
    $cellUpdates = $updates[$sheetName]
    foreach ($cellRef in $cellUpdates.Keys) {
        $cellRefStr = [string]$cellRef
        $ws.Range($cellRefStr).Value = $cellUpdates[$cellRef]
    }
}
